# Append a new entry to the end of the testing log table (first table in
# the document): Date=5/04, Test type=Expected, What I am testing=
# "Added comments and constats", Expected Outcome="The code should run
# without change", Outcome="The program ran without change",
# Changes Made="N/A".

$d = $word.ActiveDocument

$tbl = $d.Tables.Item(1)
$newRow = $tbl.Rows.Add()
$idx = $newRow.Index

$tbl.Cell($idx, 1).Range.Text = "5/04"
$tbl.Cell($idx, 2).Range.Text = "Expected"
$tbl.Cell($idx, 3).Range.Text = "Added comments and constats"
$tbl.Cell($idx, 4).Range.Text = "The code should run without change"
$tbl.Cell($idx, 5).Range.Text = "The program ran without change"
$tbl.Cell($idx, 6).Range.Text = "N/A"
